# Update the footer "Date" placeholder text from 3/4/23 to 3/22/23
# across the slide master, every slide layout, the handout master and
# the notes master (the Header/Footer "Fixed" date shown on slides).

$oldDate = "3/4/23"
$newDate = "3/22/23"

function Update-DateTextInShapes {
    param($shapes, [string]$old, [string]$new)

    $updated = 0
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tf = $shp.TextFrame
            if ($tf.HasText -and $tf.TextRange.Text -eq $old) {
                $tf.TextRange.Text = $new
                $updated = $updated + 1
            }
        }
    }
    return $updated
}

$p = $ppt.ActivePresentation
$totalUpdated = 0

# 1) Slide master footer date placeholder.
$master = $p.SlideMaster
$totalUpdated += Update-DateTextInShapes $master.Shapes $oldDate $newDate

# 2) Every slide layout's footer date placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    $totalUpdated += Update-DateTextInShapes $layout.Shapes $oldDate $newDate
}

# 3) Handout master date placeholder. Directly setting the shape's
#    TextRange is not reflected back into the handout master, so go
#    through the HeadersFooters.DateAndTime API instead.
$handoutMaster = $p.HandoutMaster
$hmDateShape = $null
for ($i = 1; $i -le $handoutMaster.Shapes.Count; $i++) {
    $shp = $handoutMaster.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq $oldDate) {
        $hmDateShape = $shp
    }
}
if ($hmDateShape -ne $null) {
    $handoutMaster.HeadersFooters.DateAndTime.Text = $newDate
    $totalUpdated += 1
}

# 4) Notes master date placeholder (same caveat as the handout master).
$notesMaster = $p.NotesMaster
$nmDateShape = $null
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shp = $notesMaster.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq $oldDate) {
        $nmDateShape = $shp
    }
}
if ($nmDateShape -ne $null) {
    $notesMaster.HeadersFooters.DateAndTime.Text = $newDate
    $totalUpdated += 1
}

Write-Host "Updated $totalUpdated date placeholder(s) from $oldDate to $newDate"
